# Weekly price-list refresh for Fruta, Feria Lagunitas de Puerto Montt - Granada
# (rows 6-39 reshuffled to the new week's order; a couple of weighted-average
# prices corrected along the way)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("D6").Value = 44305
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 17500
$ws.Range("O6").Value = 17500
$ws.Range("P6").Value = 17500
$ws.Range("S6").Value = 1167

# Row 7
$ws.Range("D7").Value = 44305
$ws.Range("M7").Value = 120
$ws.Range("O7").Value = 15000
$ws.Range("S7").Value = 967

# Row 8
$ws.Range("D8").Value = 44293
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 14000
$ws.Range("P8").Value = 14500
$ws.Range("S8").Value = 967

# Row 9
$ws.Range("D9").Value = 44295
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 160
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("S9").Value = 967

# Row 10
$ws.Range("D10").Value = 44336
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 17000
$ws.Range("O10").Value = 17000
$ws.Range("P10").Value = 17000
$ws.Range("S10").Value = 1133

# Row 11
$ws.Range("D11").Value = 44336
$ws.Range("M11").Value = 120
$ws.Range("O11").Value = 14500
$ws.Range("P11").Value = 14250
$ws.Range("S11").Value = 950

# Row 12
$ws.Range("D12").Value = 44351

# Row 13
$ws.Range("D13").Value = 44351

# Row 14
$ws.Range("D14").Value = 44348
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("S14").Value = 1000

# Row 15
$ws.Range("D15").Value = 44348
$ws.Range("N15").Value = 13000
$ws.Range("O15").Value = 13500
$ws.Range("P15").Value = 13250
$ws.Range("S15").Value = 883

# Row 16
$ws.Range("D16").Value = 44309
$ws.Range("N16").Value = 17500
$ws.Range("O16").Value = 17500
$ws.Range("P16").Value = 17500
$ws.Range("S16").Value = 1167

# Row 17
$ws.Range("D17").Value = 44309
$ws.Range("M17").Value = 200
$ws.Range("O17").Value = 14500
$ws.Range("P17").Value = 14250
$ws.Range("S17").Value = 950

# Row 18
$ws.Range("D18").Value = 44292
$ws.Range("L18").Value = "Segunda"

# Row 19
$ws.Range("D19").Value = 44299
$ws.Range("M19").Value = 60
$ws.Range("N19").Value = 17500
$ws.Range("O19").Value = 17500
$ws.Range("P19").Value = 17500
$ws.Range("S19").Value = 1167

# Row 20
$ws.Range("D20").Value = 44299
$ws.Range("M20").Value = 120
$ws.Range("O20").Value = 15000
$ws.Range("S20").Value = 967

# Row 21
$ws.Range("D21").Value = 44327
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 17000
$ws.Range("O21").Value = 17000
$ws.Range("P21").Value = 17000
$ws.Range("S21").Value = 1133

# Row 22
$ws.Range("D22").Value = 44327
$ws.Range("M22").Value = 200
$ws.Range("O22").Value = 14500
$ws.Range("P22").Value = 14250
$ws.Range("S22").Value = 950

# Row 23
$ws.Range("D23").Value = 44316
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 17500
$ws.Range("O23").Value = 17500
$ws.Range("P23").Value = 17500
$ws.Range("S23").Value = 1167

# Row 24
$ws.Range("D24").Value = 44316
$ws.Range("L24").Value = "Segunda"
$ws.Range("M24").Value = 200
$ws.Range("N24").Value = 14000
$ws.Range("O24").Value = 14500
$ws.Range("P24").Value = 14250
$ws.Range("S24").Value = 950

# Row 25
$ws.Range("D25").Value = 44330
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 17000
$ws.Range("O25").Value = 17000
$ws.Range("P25").Value = 17000
$ws.Range("S25").Value = 1133

# Row 26
$ws.Range("D26").Value = 44330
$ws.Range("M26").Value = 200
$ws.Range("O26").Value = 14500
$ws.Range("P26").Value = 14250
$ws.Range("S26").Value = 950

# Row 27
$ws.Range("D27").Value = 44323
$ws.Range("N27").Value = 17000
$ws.Range("O27").Value = 17000
$ws.Range("P27").Value = 17000
$ws.Range("S27").Value = 1133

# Row 28
$ws.Range("D28").Value = 44323
$ws.Range("M28").Value = 100
$ws.Range("O28").Value = 14000
$ws.Range("P28").Value = 14000
$ws.Range("S28").Value = 933

# Row 29
$ws.Range("D29").Value = 44306
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 100
$ws.Range("N29").Value = 17500
$ws.Range("O29").Value = 17500
$ws.Range("P29").Value = 17500
$ws.Range("S29").Value = 1167

# Row 30
$ws.Range("D30").Value = 44306
$ws.Range("L30").Value = "Segunda"
$ws.Range("M30").Value = 200
$ws.Range("N30").Value = 14000
$ws.Range("O30").Value = 14500
$ws.Range("P30").Value = 14250
$ws.Range("Q30").Value = "$/caja 15 kilos empedrada"
$ws.Range("S30").Value = 950
$ws.Range("T30").Value = 15

# Row 31
$ws.Range("D31").Value = 44285
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 160
$ws.Range("N31").Value = 15000
$ws.Range("P31").Value = 15500
$ws.Range("Q31").Value = "$/caja 15 kilos empedrada"
$ws.Range("S31").Value = 1033
$ws.Range("T31").Value = 15

# Row 32
$ws.Range("D32").Value = 44298
$ws.Range("O32").Value = 15000
$ws.Range("P32").Value = 14500
$ws.Range("Q32").Value = "$/caja 15 kilos empedrada"
$ws.Range("S32").Value = 967
$ws.Range("T32").Value = 15

# Row 33
$ws.Range("D33").Value = 44302
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 100
$ws.Range("N33").Value = 17500
$ws.Range("O33").Value = 17500
$ws.Range("P33").Value = 17500
$ws.Range("S33").Value = 1167

# Row 34
$ws.Range("D34").Value = 44302
$ws.Range("L34").Value = "Segunda"
$ws.Range("M34").Value = 200
$ws.Range("N34").Value = 14000
$ws.Range("O34").Value = 15000
$ws.Range("P34").Value = 14500
$ws.Range("S34").Value = 967

# Row 35
$ws.Range("D35").Value = 44344
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 100
$ws.Range("N35").Value = 16000
$ws.Range("O35").Value = 16000
$ws.Range("P35").Value = 16000
$ws.Range("S35").Value = 1067

# Row 36
$ws.Range("D36").Value = 44344
$ws.Range("L36").Value = "Segunda"
$ws.Range("M36").Value = 120
$ws.Range("N36").Value = 13000
$ws.Range("O36").Value = 13500
$ws.Range("P36").Value = 13250
$ws.Range("S36").Value = 883

# Row 37
$ws.Range("D37").Value = 44313
$ws.Range("L37").Value = "Especial"
$ws.Range("M37").Value = 100
$ws.Range("N37").Value = 17500
$ws.Range("O37").Value = 17500
$ws.Range("P37").Value = 17500
$ws.Range("Q37").Value = "$/caja 14 kilos empedrada"
$ws.Range("S37").Value = 1250
$ws.Range("T37").Value = 14

# Row 38
$ws.Range("D38").Value = 44313
$ws.Range("Q38").Value = "$/caja 14 kilos empedrada"
$ws.Range("S38").Value = 1143
$ws.Range("T38").Value = 14

# Row 39
$ws.Range("D39").Value = 44313
$ws.Range("M39").Value = 80
$ws.Range("N39").Value = 14000
$ws.Range("O39").Value = 14000
$ws.Range("P39").Value = 14000
$ws.Range("Q39").Value = "$/caja 14 kilos empedrada"
$ws.Range("S39").Value = 1000
$ws.Range("T39").Value = 14
